$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (all columns A:AH)
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45064.50694444445
$data[0,1] = 20.178
$data[0,2] = 13.652
$data[0,3] = 4.068
$data[0,4] = 42.752
$data[0,5] = 34.691
$data[0,6] = 15.879
$data[0,7] = 51.007
$data[0,8] = 24.432
$data[0,9] = 10.251
$data[0,10] = 15.607
$data[0,11] = 16.868
$data[0,12] = 17.597
$data[0,13] = 5.069
$data[0,14] = 15.79
$data[0,15] = 22.094
$data[0,16] = 13.41
$data[0,17] = 3.498
$data[0,18] = 2.451
$data[0,19] = 232.996
$data[0,20] = 43.923
$data[0,21] = 14.575
$data[0,22] = 29.082
$data[0,23] = 14.996
$data[0,24] = 3.13
$data[0,25] = 25.159
$data[0,26] = 12.874
$data[0,27] = 11.654
$data[0,28] = 13.651
$data[0,29] = 17.331
$data[0,30] = 3.457
$data[0,31] = 45.214
$data[0,32] = 8.105
$data[0,33] = 18.222
$data[1,0] = 45064.51388888889
$data[1,1] = 24.021
$data[1,2] = 17.339
$data[1,3] = 2.064
$data[1,4] = 51.919
$data[1,5] = 42.489
$data[1,6] = 18.903
$data[1,7] = 72.117
$data[1,8] = 29.086
$data[1,9] = 12.736
$data[1,10] = 18.975
$data[1,11] = 20.81
$data[1,12] = 21.827
$data[1,13] = 6.038
$data[1,14] = 18.798
$data[1,15] = 26.631
$data[1,16] = 15.959
$data[1,17] = 1.631
$data[1,18] = 1.31
$data[1,19] = 278.82
$data[1,20] = 52.574
$data[1,21] = 17.351
$data[1,22] = 35.174
$data[1,23] = 18.451
$data[1,24] = 3.055
$data[1,25] = 34.987
$data[1,26] = 15.326
$data[1,27] = 13.72
$data[1,28] = 16.102
$data[1,29] = 21.639
$data[1,30] = 1.266
$data[1,31] = 65.511
$data[1,32] = 9.738
$data[1,33] = 21.692
$data[2,0] = 45064.52083333334
$data[2,1] = 6.246
$data[2,2] = 4.231
$data[2,3] = 0.98
$data[2,4] = 13.411
$data[2,5] = 10.739
$data[2,6] = 4.916
$data[2,7] = 24.951
$data[2,8] = 7.562
$data[2,9] = 3.225
$data[2,10] = 4.655
$data[2,11] = 5.397
$data[2,12] = 5.566
$data[2,13] = 1.578
$data[2,14] = 4.888
$data[2,15] = 6.88
$data[2,16] = 4.363
$data[2,17] = 0.972
$data[2,18] = 0.524
$data[2,19] = 67.11499999999999
$data[2,20] = 13.948
$data[2,21] = 4.511
$data[2,22] = 9.154999999999999
$data[2,23] = 4.758
$data[2,24] = 1.049
$data[2,25] = 11.214
$data[2,26] = 3.985
$data[2,27] = 3.696
$data[2,28] = 4.314
$data[2,29] = 5.519
$data[2,30] = 0.773
$data[2,31] = 23.009
$data[2,32] = 2.445
$data[2,33] = 5.642
$data[3,0] = 45064.52777777778
$data[3,1] = 2.4
$data[3,2] = 1.46
$data[3,3] = 0.64
$data[3,4] = 5.12
$data[3,5] = 3.93
$data[3,6] = 1.89
$data[3,7] = 10.97
$data[3,8] = 2.91
$data[3,9] = 1.19
$data[3,10] = 1.61
$data[3,11] = 2.08
$data[3,12] = 2.08
$data[3,13] = 0.61
$data[3,14] = 1.88
$data[3,15] = 2.62
$data[3,16] = 1.81
$data[3,17] = 0.7
$data[3,18] = 0.3
$data[3,19] = 21.32
$data[3,20] = 5.47
$data[3,21] = 1.74
$data[3,22] = 3.5
$data[3,23] = 1.82
$data[3,24] = 0.5600000000000001
$data[3,25] = 4.72
$data[3,26] = 1.53
$data[3,27] = 1.5
$data[3,28] = 1.74
$data[3,29] = 2.07
$data[3,30] = 0.5600000000000001
$data[3,31] = 10.12
$data[3,32] = 0.88
$data[3,33] = 2.17
$ws.Range("A2:AH5").Value = $data

# Delete row 6 entirely (reduces used range to A1:AH5)
$ws.Range("A6:AH6").Delete()

# Adjust column widths to match target widths
$ws.Range("B1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("C1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("G1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("J1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("K1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("L1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("M1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("O1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("P1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("T1").EntireColumn.ColumnWidth = 8.166666666666666
$ws.Range("V1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Z1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AA1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AB1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AC1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AD1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AH1").EntireColumn.ColumnWidth = 7.166666666666667
